$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1 "Play Book of Santa Slot Game for Free - Review").
# ---------------------------------------------------------------------
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd() -eq "Play Book of Santa Slot Game for Free - Review" -and `
        $cand.Style.NameLocal -eq "Heading 1") {
        $titlePara = $cand
        break
    }
}

$titlePara.Range.InsertParagraphAfter()
$metaPara = $titlePara.Next()
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaStart = $metaRange.Start

$boldText = "Meta description"
$restText = ": Read our review of Book of Santa slot game, play it for free and enter the Christmas atmosphere with expandable symbols and free spin feature."

$metaRange.InsertBefore($boldText + $restText)

$boldRange = $d.Range($metaStart, $metaStart + $boldText.Length)
$boldRange.Bold = 1

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the
#    document ("Play Book of Santa Slot Game for Free - Review", the
#    Normal/bold copy -- not the Heading1 title) and
# 3) replace the text of the italic paragraph that follows it with
#    the new feature-image prompt, keeping its italic formatting.
# ---------------------------------------------------------------------
$dupTitlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd() -eq "Play Book of Santa Slot Game for Free - Review" -and `
        $cand.Style.NameLocal -ne "Heading 1") {
        $dupTitlePara = $cand
        break
    }
}

$dupTitlePara.Range.Delete()

$oldBlurb = "Read our review of Book of Santa slot game, play it for free and enter the Christmas atmosphere with expandable symbols and free spin feature."
$blurbPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq $oldBlurb) {
        $blurbPara = $d.Paragraphs($i)
        break
    }
}

$blurbRange = $blurbPara.Range
$blurbTextRange = $d.Range($blurbRange.Start, $blurbRange.End - 1)
$blurbTextRange.Text = "Create a feature image for Book of Santa that captures the essence of the game's theme and unique character. The image should be in a cartoon style and feature a happy-looking Maya warrior wearing glasses, as the protagonist of the game. The Maya warrior should be holding a large book in his hands, which should resemble Santa's book, with a few gifts spilling out of it. The background of the image should be a cozy fireplace scene, with the grid of the slot game superimposed on it. The image should use bright and cheerful colors and convey the festive mood of the holiday season."
